# Update Name of Algo
# Apply updated numeric values produced by a re-run of the KNN imputation
# algorithm to the relevant cells on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -10.872
$ws.Range("D3").Value = -7.136000000000001
$ws.Range("B12").Value = 4.806
$ws.Range("C14").Value = -12.21
$ws.Range("D20").Value = -7.755000000000001
$ws.Range("D25").Value = -8.163
$ws.Range("C26").Value = -13.088
$ws.Range("B27").Value = 5.441000000000001
$ws.Range("D30").Value = -7.176
$ws.Range("C31").Value = -12.319
$ws.Range("B32").Value = 6.205
$ws.Range("C35").Value = -12.35
$ws.Range("B36").Value = 8.673999999999999
$ws.Range("C37").Value = -13.665
$ws.Range("B38").Value = 5.148000000000001
$ws.Range("D44").Value = -7.749
$ws.Range("C45").Value = -12.776
$ws.Range("B46").Value = 5.839
$ws.Range("D47").Value = -7.375
$ws.Range("C52").Value = -11.07
$ws.Range("B54").Value = 5.006
$ws.Range("B55").Value = 4.742
$ws.Range("B56").Value = 4.49
$ws.Range("C57").Value = -13.697
$ws.Range("D58").Value = -8.021000000000001
$ws.Range("B67").Value = 5.169
$ws.Range("B69").Value = 5.146999999999999
$ws.Range("B72").Value = 5.501
$ws.Range("D78").Value = -7.906999999999999
$ws.Range("C81").Value = -13.319
$ws.Range("B83").Value = 5.412000000000001
$ws.Range("C83").Value = -13.772
$ws.Range("D84").Value = -8.360000000000001
$ws.Range("B86").Value = 5.052
$ws.Range("D89").Value = -6.856999999999999
$ws.Range("B91").Value = 5.949
$ws.Range("D91").Value = -6.395
$ws.Range("D92").Value = -6.672
$ws.Range("B93").Value = 5.685999999999999
$ws.Range("D96").Value = -7.38
$ws.Range("B99").Value = 5.718
$ws.Range("C100").Value = -12.807
$ws.Range("C102").Value = -13.583
$ws.Range("D102").Value = -7.723999999999999
